$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 153 (existing rows 153:164 shift down to 155:166).
$ws.Range("A153:A154").EntireRow.Insert()

# New row 153 (Primera) - new week's data.
$ws.Range("A153").Value = 11
$ws.Range("B153").Value = "Vega Monumental Concepción"
$ws.Range("C153").Value = "Bíobío"
$ws.Range("D153").Value = 44491
$ws.Range("E153").Value = 8
$ws.Range("F153").Value = 100112023
$ws.Range("G153").Value = "Brócoli"
$ws.Range("H153").Value = "Sin especificar"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 1000
$ws.Range("K153").Value = 800
$ws.Range("L153").Value = 900
$ws.Range("M153").Value = 850
$ws.Range("N153").Value = "`$/unidad"
$ws.Range("O153").Value = "Región Metropolitana"
$ws.Range("P153").Value = 850
$ws.Range("Q153").Value = 1
$ws.Range("R153").Value = "Hortaliza"

# New row 154 (Segunda) - new week's data.
$ws.Range("A154").Value = 11
$ws.Range("B154").Value = "Vega Monumental Concepción"
$ws.Range("C154").Value = "Bíobío"
$ws.Range("D154").Value = 44491
$ws.Range("E154").Value = 8
$ws.Range("F154").Value = 100112023
$ws.Range("G154").Value = "Brócoli"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Segunda"
$ws.Range("J154").Value = 500
$ws.Range("K154").Value = 700
$ws.Range("L154").Value = 700
$ws.Range("M154").Value = 700
$ws.Range("N154").Value = "`$/unidad"
$ws.Range("O154").Value = "Región Metropolitana"
$ws.Range("P154").Value = 700
$ws.Range("Q154").Value = 1
$ws.Range("R154").Value = "Hortaliza"

Write-Host "Inserted new weekly rows 153-154; dimension now $($ws.Range("A1").CurrentRegion.Address)"
